$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing O1 header cell onto the new P1:R1
# header cells so they pick up the same style (bold, centered, bordered).
$ws.Range("O1").Copy()
$ws.Range("P1:R1").PasteSpecial(-4122)

# New header cells in row 1 (P, Q, R)
$ws.Range("P1").Value = "(부)문화미디어 전기"
$ws.Range("Q1").Value = "(부)문화미디어 전필"
$ws.Range("R1").Value = "(부)문화미디어 전선"

# Data rows
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 15

$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0

$ws.Range("P4").Value = 6
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 15
